$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "DataSet" column (column D) entirely - this shifts columns E:M left to D:L
$ws.Range("D1").EntireColumn.Delete()

# Rename "RECALL" header (now in column K) to "Average RECALL"
$ws.Range("K1").Value = "Average RECALL"

# Update the active selection to match the post-edit state
$ws.Range("G10").Select()
